$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.191.60"
$ws.Range("E2").Value = "  -3.11%  "
$ws.Range("D3").Value = "'2.518.81"
$ws.Range("E3").Value = "  -5.47%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'574.22"
$ws.Range("E5").Value = "  -3.99%  "
$ws.Range("D6").Value = "'169.78"
$ws.Range("E6").Value = "  -3.25%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.510"
$ws.Range("E8").Value = "  -2.83%  "
$ws.Range("D9").Value = "'2.518.78"
$ws.Range("E9").Value = "  -5.54%  "
$ws.Range("D10").Value = "'0.159"
$ws.Range("E10").Value = "  -5.96%  "
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").Value = "'0.343"
$ws.Range("E12").Value = "  -4.08%  "
$ws.Range("D13").Value = "'4.82"
$ws.Range("E13").Value = "  -3.75%  "
$ws.Range("D14").Value = "'2.980.77"
$ws.Range("E14").Value = "  -5.63%  "
$ws.Range("D15").Value = "'70.001.51"
$ws.Range("E15").Value = "  -3.29%  "
$ws.Range("D16").Value = "'0.0000179"
$ws.Range("E16").Value = "  -3.54%  "
$ws.Range("D17").Value = "'24.95"
$ws.Range("E17").Value = "  -5.45%  "
$ws.Range("D18").Value = "'2.514.81"
$ws.Range("E18").Value = "  -5.78%  "
$ws.Range("D19").Value = "'11.51"
$ws.Range("E19").Value = "  -5.99%  "
$ws.Range("D20").Value = "'7.57"
$ws.Range("E20").Value = "  -8.57%  "
$ws.Range("D21").Value = "'353.98"
$ws.Range("E21").Value = "  -4.55%  "
$ws.Range("D22").Value = "'3.94"
$ws.Range("E22").Value = "  -5.99%  "
$ws.Range("D23").Value = "'1.98"
$ws.Range("E23").Value = "  -3.04%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "'68.97"
$ws.Range("E25").Value = "  -4.40%  "
$ws.Range("D26").Value = "'4.09"
$ws.Range("E26").Value = "  -5.80%  "
$ws.Range("E27").Value = "  -5.93%  "
$ws.Range("D28").Value = "'2.649.12"
$ws.Range("E28").Value = "  -5.65%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "'0.0₃0909"
$ws.Range("E30").Value = "  -6.59%  "
$ws.Range("D31").Value = "'7.83"
$ws.Range("E31").Value = "  -3.59%  "
$ws.Range("D32").Value = "'480.81"
$ws.Range("E32").Value = "  -4.22%  "
$ws.Range("E33").Value = "  -0.95%  "
$ws.Range("D34").Value = "'1.76"
$ws.Range("E34").Value = "  -3.84%  "
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("E36").Value = "  +4.33%  "
$ws.Range("D37").Value = "'156.18"
$ws.Range("E37").Value = "  -4.27%  "
$ws.Range("D38").Value = "'18.87"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").Value = "'18.53"
$ws.Range("E39").Value = "  -5.22%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").Value = "'1.23"
$ws.Range("E41").Value = "  -11.44%  "
$ws.Range("E42").Value = "  -7.64%  "
$ws.Range("D43").Value = "'0.320"
$ws.Range("E43").Value = "  -3.97%  "
$ws.Range("D44").Value = "'4.73"
$ws.Range("E44").Value = "  -5.52%  "
$ws.Range("E45").Value = "  -7.50%  "
$ws.Range("D46").Value = "'38.29"
$ws.Range("E46").Value = "  -3.14%  "
$ws.Range("D47").Value = "'141.98"
$ws.Range("E47").Value = "  -9.27%  "
$ws.Range("D48").Value = "'3.52"
$ws.Range("E48").Value = "  -6.31%  "
$ws.Range("D49").Value = "'0.526"
$ws.Range("E49").Value = "  -6.04%  "
$ws.Range("E50").Value = "  -6.90%  "
$ws.Range("E51").Value = "  -1.41%  "
